$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format first so numeric-looking strings
# (e.g. "581.51") are stored as text, matching the source data, not numbers.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "64.121.83"
$ws.Range("E2").Value = "  +1.99%  "
$ws.Range("D3").Value = "2.533.12"
$ws.Range("E3").Value = "  +2.80%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "581.51"
$ws.Range("E5").Value = "  +1.09%  "
$ws.Range("D6").Value = "151.96"
$ws.Range("E6").Value = "  +4.09%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "0.538"
$ws.Range("E8").Value = "  +0.30%  "
$ws.Range("D9").Value = "0.112"
$ws.Range("E9").Value = "  +0.74%  "
$ws.Range("E10").Value = "  -0.39%  "
$ws.Range("D11").Value = "5.28"
$ws.Range("E11").Value = "  +0.06%  "
$ws.Range("D12").Value = "0.355"
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("D13").Value = "29.47"
$ws.Range("E13").Value = "  +1.40%  "
$ws.Range("D14").Value = "0.0000179"
$ws.Range("E14").Value = "  +1.06%  "
$ws.Range("D15").Value = "2.996.47"
$ws.Range("E15").Value = "  +2.94%  "
$ws.Range("D16").Value = "63.560.10"
$ws.Range("E16").Value = "  +1.25%  "
$ws.Range("D17").Value = "2.540.51"
$ws.Range("E17").Value = "  +3.05%  "
$ws.Range("D18").Value = "7.85"
$ws.Range("E18").Value = "  -1.74%  "
$ws.Range("D19").Value = "10.97"
$ws.Range("E19").Value = "  -0.44%  "
$ws.Range("D20").Value = "4.26"
$ws.Range("E20").Value = "  +3.02%  "
$ws.Range("D21").Value = "327.50"
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("E22").Value = "  +0.71%  "
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").Value = "10.10"
$ws.Range("E24").Value = "  -1.20%  "
$ws.Range("D25").Value = "65.56"
$ws.Range("E25").Value = "  -0.49%  "
$ws.Range("D26").Value = "658.91"
$ws.Range("E26").Value = "  +1.10%  "
$ws.Range("D27").Value = "0.0000103"
$ws.Range("E27").Value = "  +5.17%  "
$ws.Range("D28").Value = "2.670.80"
$ws.Range("E28").Value = "  +3.36%  "
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("D30").Value = "1.48"
$ws.Range("E30").Value = "  +1.84%  "
$ws.Range("D31").Value = "8.06"
$ws.Range("E31").Value = "  +0.68%  "
$ws.Range("D32").Value = "1.86"
$ws.Range("E32").Value = "  +0.27%  "
$ws.Range("D33").Value = "0.135"
$ws.Range("E33").Value = "  +1.32%  "
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("D35").Value = "1.54"
$ws.Range("E35").Value = "  +1.11%  "
$ws.Range("D36").Value = "4.82"
$ws.Range("E36").Value = "  +1.31%  "
$ws.Range("D37").Value = "5.55"
$ws.Range("E37").Value = "  +2.93%  "
$ws.Range("D38").Value = "2.82"
$ws.Range("E38").Value = "  +2.25%  "
$ws.Range("D39").Value = "0.371"
$ws.Range("E39").Value = "  +0.47%  "
$ws.Range("D40").Value = "18.87"
$ws.Range("E40").Value = "  +0.79%  "
$ws.Range("D41").Value = "151.51"
$ws.Range("E41").Value = "  +0.32%  "
$ws.Range("E42").Value = "  +2.77%  "
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").Value = "158.81"
$ws.Range("E44").Value = "  +3.12%  "
$ws.Range("E45").Value = "  -4.37%  "
$ws.Range("D46").Value = "15.44"
$ws.Range("E46").Value = "  +1.35%  "
$ws.Range("D47").Value = "3.64"
$ws.Range("E47").Value = "  +1.67%  "
$ws.Range("D48").Value = "20.92"
$ws.Range("E48").Value = "  +2.67%  "
$ws.Range("D49").Value = "0.620"
$ws.Range("E49").Value = "  +2.25%  "
$ws.Range("D50").Value = "0.0520"
$ws.Range("E50").Value = "  +1.49%  "
$ws.Range("D51").Value = "0.0229"
$ws.Range("E51").Value = "  +1.87%  "

# Restore the default (Normal) style on column D so no residual
# number-format styling is left behind, matching the original file.
$dRange.Style = "Normal"
